# This script applies the "dynamic start and end times" bug fix to the
# Schedule worksheet: the earliest four 15-minute time slots (11:00-11:45)
# are dropped and a new slot (19:15) already existed at the far end, so the
# whole time axis effectively shifts one hour later and the four now-unused
# trailing columns are removed. The per-row rotation numbers are
# recalculated to line up with the new column positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 0) The gray "last slot" highlight format (currently on AG2 and AH6) needs to
#    land on the new final populated cell of those rows (AC2 and AD6) after
#    the shift. Grab a copy of that exact cell format now, before the source
#    columns are removed, so the existing style entry is reused rather than a
#    new one being created.
$ws.Range("AG2").Copy() | Out-Null
$ws.Range("AC2").PasteSpecial(-4122) | Out-Null
$ws.Range("AH6").Copy() | Out-Null
$ws.Range("AD6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 1) Remove the four trailing time-slot columns (AF:AI) that are no longer
#    needed, shrinking the sheet from A1:AI12 down to A1:AE12.
$ws.Range("AF1:AI12").EntireColumn.Delete() | Out-Null

# 2) Re-label the header row with the shifted time slots.
$ws.Range("B1").Value = "12:00"
$ws.Range("C1").Value = "12:15"
$ws.Range("D1").Value = "12:30"
$ws.Range("E1").Value = "12:45"
$ws.Range("F1").Value = "13:00"
$ws.Range("G1").Value = "13:15"
$ws.Range("H1").Value = "13:30"
$ws.Range("I1").Value = "13:45"
$ws.Range("J1").Value = "14:00"
$ws.Range("K1").Value = "14:15"
$ws.Range("L1").Value = "14:30"
$ws.Range("M1").Value = "14:45"
$ws.Range("N1").Value = "15:00"
$ws.Range("O1").Value = "15:15"
$ws.Range("P1").Value = "15:30"
$ws.Range("Q1").Value = "15:45"
$ws.Range("R1").Value = "16:00"
$ws.Range("S1").Value = "16:15"
$ws.Range("T1").Value = "16:30"
$ws.Range("U1").Value = "16:45"
$ws.Range("V1").Value = "17:00"
$ws.Range("W1").Value = "17:15"
$ws.Range("X1").Value = "17:30"
$ws.Range("Y1").Value = "17:45"
$ws.Range("Z1").Value = "18:00"
$ws.Range("AA1").Value = "18:15"
$ws.Range("AB1").Value = "18:30"
$ws.Range("AC1").Value = "18:45"
$ws.Range("AD1").Value = "19:00"
$ws.Range("AE1").Value = "19:15"

# 3) Recompute the rotation numbers for every attraction row so each value
#    sits under its correct (now shifted) time column.
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 5
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4
$ws.Range("N2").Value = 3
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 6
$ws.Range("Q2").Value = 5
$ws.Range("R2").Value = 7
$ws.Range("S2").Value = 4
$ws.Range("T2").Value = 3
$ws.Range("U2").Value = 2
$ws.Range("V2").Value = 6
$ws.Range("W2").Value = 5
$ws.Range("X2").Value = 7
$ws.Range("Y2").Value = 4
$ws.Range("Z2").Value = 3
$ws.Range("AA2").Value = 2
$ws.Range("AB2").Value = 6
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 2
$ws.Range("I3").Value = 3
$ws.Range("J3").Value = 6
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 5
$ws.Range("M3").Value = 1
$ws.Range("N3").Value = 4
$ws.Range("O3").Value = 3
$ws.Range("P3").Value = 2
$ws.Range("Q3").Value = 6
$ws.Range("R3").Value = 5
$ws.Range("S3").Value = 7
$ws.Range("T3").Value = 4
$ws.Range("U3").Value = 3
$ws.Range("V3").Value = 2
$ws.Range("W3").Value = 6
$ws.Range("X3").Value = 5
$ws.Range("Y3").Value = 7
$ws.Range("Z3").Value = 4
$ws.Range("AA3").Value = 3
$ws.Range("AB3").Value = 2
$ws.Range("AC3").Value = 6
$ws.Range("AE3").Value = 4
$ws.Range("F4").Value = 3
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 2
$ws.Range("K4").Value = 6
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 5
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 4
$ws.Range("P4").Value = 3
$ws.Range("Q4").Value = 2
$ws.Range("R4").Value = 6
$ws.Range("S4").Value = 5
$ws.Range("T4").Value = 7
$ws.Range("U4").Value = 4
$ws.Range("V4").Value = 3
$ws.Range("W4").Value = 2
$ws.Range("X4").Value = 6
$ws.Range("Y4").Value = 5
$ws.Range("Z4").Value = 7
$ws.Range("AA4").Value = 4
$ws.Range("AB4").Value = 3
$ws.Range("AC4").Value = 2
$ws.Range("AD4").Value = 6
$ws.Range("G5").Value = 3
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 4
$ws.Range("L5").Value = 6
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 5
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 4
$ws.Range("Q5").Value = 3
$ws.Range("R5").Value = 2
$ws.Range("S5").Value = 6
$ws.Range("T5").Value = 5
$ws.Range("U5").Value = 7
$ws.Range("V5").Value = 4
$ws.Range("W5").Value = 3
$ws.Range("X5").Value = 2
$ws.Range("Y5").Value = 6
$ws.Range("Z5").Value = 5
$ws.Range("AA5").Value = 7
$ws.Range("AB5").Value = 4
$ws.Range("AC5").Value = 3
$ws.Range("AD5").Value = 2
$ws.Range("AE5").Value = 6
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 4
$ws.Range("L6").Value = 3
$ws.Range("M6").Value = 6
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 5
$ws.Range("P6").Value = 7
$ws.Range("Q6").Value = 4
$ws.Range("R6").Value = 3
$ws.Range("S6").Value = 2
$ws.Range("T6").Value = 6
$ws.Range("U6").Value = 5
$ws.Range("V6").Value = 7
$ws.Range("W6").Value = 4
$ws.Range("X6").Value = 3
$ws.Range("Y6").Value = 2
$ws.Range("Z6").Value = 6
$ws.Range("AA6").Value = 5
$ws.Range("AB6").Value = 7
$ws.Range("AC6").Value = 4
$ws.Range("AE6").Value = 2
$ws.Range("N7").Value = 6
$ws.Range("O7").Value = 0
$ws.Range("H10").Value = 3
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 5
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 4
$ws.Range("M10").Value = 3
$ws.Range("N10").Value = 2
$ws.Range("O10").Value = 6
$ws.Range("P10").Value = 5
$ws.Range("Q10").Value = 7
$ws.Range("R10").Value = 4
$ws.Range("S10").Value = 3
$ws.Range("T10").Value = 2
$ws.Range("U10").Value = 6
$ws.Range("V10").Value = 5
$ws.Range("W10").Value = 7
$ws.Range("X10").Value = 4
$ws.Range("Y10").Value = 3
$ws.Range("Z10").Value = 2
$ws.Range("AA10").Value = 6
$ws.Range("AB10").Value = 5
$ws.Range("AC10").Value = 7
$ws.Range("AD10").Value = 4
$ws.Range("AE10").Value = 5

# 4) A few cells that used to hold a rotation number no longer apply under
#    the new time range, so clear them back out.
$ws.Range("AD2").ClearContents() | Out-Null
$ws.Range("AE2").ClearContents() | Out-Null
$ws.Range("R7").ClearContents() | Out-Null
$ws.Range("S7").ClearContents() | Out-Null

# 5) The light-gray highlight format was already copied onto AC2/AD6 in step
#    0; just set their final values here.
$ws.Range("AC2").Value = 5
$ws.Range("AD6").Value = 5
